# Applies updated market-price / profit figures to each Leve sheet,
# mirroring the scheduled price-refresh runner output.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 288.42856
$ws.Range("I12").Value = 303.66666
$ws.Range("K12").Value = 303.66666
$ws.Range("M12").Value = -133.66666
$ws.Range("H16").Value = 15555
$ws.Range("J16").Value = 15555
$ws.Range("L16").Value = 15555
$ws.Range("N16").Value = -16015
$ws.Range("H17").Value = 5575.7827
$ws.Range("J17").Value = 6011.5713
$ws.Range("L17").Value = 18034.7139
$ws.Range("N17").Value = -18370.7139
$ws.Range("H88").Value = 1686.5714
$ws.Range("J88").Value = 1686.5714
$ws.Range("L88").Value = 1686.5714
$ws.Range("N88").Value = -2498.5714
$ws.Range("H91").Value = 1686.5714
$ws.Range("J91").Value = 1686.5714
$ws.Range("L91").Value = 1686.5714
$ws.Range("N91").Value = -4494.5714
$ws.Range("H92").Value = 550
$ws.Range("J92").Value = 550
$ws.Range("L92").Value = 550
$ws.Range("N92").Value = -3046
$ws.Range("H99").Value = 1526.6666
$ws.Range("I99").Value = 2240.3333
$ws.Range("J99").Value = 813
$ws.Range("K99").Value = 6720.999899999999
$ws.Range("L99").Value = 2439
$ws.Range("M99").Value = -5222.999899999999
$ws.Range("N99").Value = -5435
$ws.Range("H101").Value = 25000374
$ws.Range("I101").Value = 33333500
$ws.Range("J101").Value = 999
$ws.Range("K101").Value = 100000500
$ws.Range("L101").Value = 2997
$ws.Range("M101").Value = -99998878
$ws.Range("N101").Value = -6241
$ws.Range("H112").Value = 2995.1667
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 2995.1667
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 8985.500100000001
$ws.Range("M112").Value = ""
$ws.Range("N112").Value = -11201.5001
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").Value = ""
$ws.Range("H116").Value = 5266.6665
$ws.Range("I116").Value = 5266.6665
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 5266.6665
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -1824.6665
$ws.Range("N116").Value = ""
$ws.Range("H125").Value = 115389480
$ws.Range("I125").Value = 125004170
$ws.Range("K125").Value = 1125037530
$ws.Range("M125").Value = -1125035070
$ws.Range("H132").Value = 2177.9512
$ws.Range("I132").Value = 1343.0834
$ws.Range("K132").Value = 4029.2502
$ws.Range("M132").Value = -1499.2502
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").Value = ""
$ws.Range("H138").Value = 3305.8
$ws.Range("J138").Value = 2612.4
$ws.Range("L138").Value = 7837.200000000001
$ws.Range("N138").Value = -18117.2

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 667566.7
$ws.Range("I10").Value = 667566.7
$ws.Range("K10").Value = 667566.7
$ws.Range("M10").Value = -667396.7
$ws.Range("H21").Value = 3627.5
$ws.Range("J21").Value = 5555
$ws.Range("L21").Value = 5555
$ws.Range("N21").Value = -6303
$ws.Range("H32").Value = 3778.7727
$ws.Range("I32").Value = 3778.7727
$ws.Range("K32").Value = 3778.7727
$ws.Range("M32").Value = -3491.7727
$ws.Range("H61").Value = 1608.8
$ws.Range("I61").Value = 1608.8
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1608.8
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1396.8
$ws.Range("N61").Value = ""
$ws.Range("H74").Value = 1307.5
$ws.Range("I74").Value = 1166.2
$ws.Range("K74").Value = 1166.2
$ws.Range("M74").Value = -292.2
$ws.Range("H77").Value = 1307.5
$ws.Range("I77").Value = 1166.2
$ws.Range("K77").Value = 5831
$ws.Range("M77").Value = -1463
$ws.Range("H132").Value = 1200
$ws.Range("I132").Value = 1200
$ws.Range("K132").Value = 3600
$ws.Range("M132").Value = -1070
$ws.Range("H136").Value = 1608.8
$ws.Range("I136").Value = 1608.8
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4826.4
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2276.4
$ws.Range("N136").Value = ""
$ws.Range("H138").Value = 2449964.5
$ws.Range("J138").Value = 2449964.5
$ws.Range("L138").Value = 2449964.5
$ws.Range("N138").Value = -2460244.5

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3815.2856
$ws.Range("J107").Value = 5186.6
$ws.Range("L107").Value = 5186.6
$ws.Range("N107").Value = -9026.6

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3335.5
$ws.Range("I31").Value = 5999.5
$ws.Range("K31").Value = 5999.5
$ws.Range("M31").Value = -5704.5
$ws.Range("H34").Value = 3335.5
$ws.Range("I34").Value = 5999.5
$ws.Range("K34").Value = 5999.5
$ws.Range("M34").Value = -5797.5
$ws.Range("H99").Value = 1500
$ws.Range("I99").Value = 1500
$ws.Range("K99").Value = 1500
$ws.Range("M99").Value = -2
$ws.Range("H126").Value = 1500
$ws.Range("I126").Value = 1500
$ws.Range("K126").Value = 4500
$ws.Range("M126").Value = -2030
$ws.Range("H134").Value = 677.9
$ws.Range("I134").Value = 696.65515
$ws.Range("J134").Value = 134
$ws.Range("K134").Value = 2089.96545
$ws.Range("L134").Value = 402
$ws.Range("M134").Value = 445.0345499999999
$ws.Range("N134").Value = -5472
$ws.Range("H141").Value = 35725
$ws.Range("J141").Value = 35063.285
$ws.Range("L141").Value = 35063.285
$ws.Range("N141").Value = -45423.285

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 50
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = ""
$ws.Range("H27").Value = 50
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").Value = ""

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 806.8461
$ws.Range("I22").Value = 804
$ws.Range("K22").Value = 804
$ws.Range("M22").Value = -509
$ws.Range("H27").Value = 806.8461
$ws.Range("I27").Value = 804
$ws.Range("K27").Value = 804
$ws.Range("M27").Value = -697
$ws.Range("H44").Value = 16200
$ws.Range("J44").Value = 16200
$ws.Range("L44").Value = 16200
$ws.Range("N44").Value = -17112
$ws.Range("H51").Value = 17000
$ws.Range("J51").Value = 17000
$ws.Range("L51").Value = 17000
$ws.Range("N51").Value = -17956
$ws.Range("H132").Value = 9155.286
$ws.Range("I132").Value = 10097.833
$ws.Range("K132").Value = 30293.499
$ws.Range("M132").Value = -27763.499
$ws.Range("H136").Value = 4258.7
$ws.Range("I136").Value = 4286.8887
$ws.Range("K136").Value = 12860.6661
$ws.Range("M136").Value = -10310.6661

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1004309
$ws.Range("I81").Value = 1848.1666
$ws.Range("J81").Value = 2508000.2
$ws.Range("K81").Value = 3696.3332
$ws.Range("L81").Value = 5016000.4
$ws.Range("M81").Value = -2635.3332
$ws.Range("N81").Value = -5018122.4
$ws.Range("H84").Value = 1004309
$ws.Range("I84").Value = 1848.1666
$ws.Range("J84").Value = 2508000.2
$ws.Range("K84").Value = 18481.666
$ws.Range("L84").Value = 25080002
$ws.Range("M84").Value = -13177.666
$ws.Range("N84").Value = -25090610
$ws.Range("H107").Value = 1000.35
$ws.Range("I107").Value = 1000.8
$ws.Range("K107").Value = 3002.4
$ws.Range("M107").Value = -1082.4
$ws.Range("H132").Value = 12217
$ws.Range("I132").Value = 12217
$ws.Range("K132").Value = 36651
$ws.Range("M132").Value = -34121
$ws.Range("H136").Value = 11162.625
$ws.Range("I136").Value = 11162.625
$ws.Range("K136").Value = 33487.875
$ws.Range("M136").Value = -30937.875

Write-Host "Applied scheduled price/profit updates."